# GEO-REQ.schema.docx edit:
#  - Rename first "Objet geolocalisationWrapper (geolocalisation)" heading
#    to "Objet geolocalisation".
#  - Remove the "geolocalisation" reference table (the table right after
#    that heading) entirely.
#  - Remove the now-orphaned "Type geolocalisation" heading paragraph
#    that used to describe that table.
# The following "Type resource" heading + table(s) are left untouched.

$d = $word.ActiveDocument

# 1) Rename the heading text.
$d.Content.Find.Execute("Objet geolocalisationWrapper (geolocalisation)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objet geolocalisation", 2)

# 2) Delete the first table (the "geolocalisation" field table that
#    immediately followed the heading we just renamed).
$d.Tables.Item(1).Delete()

# 3) Delete the orphaned "Type geolocalisation" heading paragraph,
#    including its paragraph mark, so the table that used to follow it
#    (now the first table) directly follows the renamed heading.
$rng = $d.Content
$rng.Find.Execute("Type geolocalisation")
$headingPara = $rng.Paragraphs.Item(1)
$startPos = $headingPara.Range.Start
$endPos = $d.Tables.Item(1).Range.Start
$d.Range($startPos, $endPos).Delete()
